$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.862.59"
$ws.Range("E2").Value = "  -2.40%  "
$ws.Range("D3").Value = "2.573.92"
$ws.Range("E3").Value = "  -0.23%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'302.28"
$ws.Range("E5").Value = "  +0.12%  "
$ws.Range("D6").Value = "'96.66"
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -2.63%  "
$ws.Range("D10").Value = "'36.22"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "'0.0810"
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("E12").Value = "  -2.42%  "
$ws.Range("E13").Value = "  +6.27%  "
$ws.Range("D14").Value = "2.529.67"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("E15").Value = "  -0.82%  "
$ws.Range("D16").Value = "'14.36"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "42.888.61"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "0.0₃0998"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "'12.92"
$ws.Range("E19").Value = "  +3.54%  "
$ws.Range("E20").Value = "  -1.10%  "
$ws.Range("D21").Value = "'72.11"
$ws.Range("E21").Value = "  -2.15%  "
$ws.Range("E22").Value = "  -4.53%  "
$ws.Range("D23").Value = "'2.94"
$ws.Range("E23").Value = "  -0.28%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -5.01%  "
$ws.Range("D25").Value = "'28.99"
$ws.Range("E25").Value = "  -1.80%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("D27").Value = "'10.29"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'37.64"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("E29").Value = "  -5.66%  "
$ws.Range("D30").Value = "'6.00"
$ws.Range("E30").Value = "  -3.12%  "
$ws.Range("D31").Value = "'154.61"
$ws.Range("E31").Value = "  +1.31%  "
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("E33").Value = "  -5.14%  "
$ws.Range("E34").Value = "  -1.99%  "
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("D36").Value = "'18.26"
$ws.Range("E36").Value = "  +8.04%  "
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("E38").Value = "  -0.83%  "
$ws.Range("D39").Value = "'23.14"
$ws.Range("E39").Value = "  -4.38%  "
$ws.Range("D40").Value = "'2.08"
$ws.Range("E40").Value = "  +27.59%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").Value = "'0.0311"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("B42").Value = "NEARProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D42").Value = "'3.41"
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "2.072.05"
$ws.Range("E44").Value = "  +1.61%  "
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "'9.20"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "'85.31"
$ws.Range("E47").Value = "  -3.19%  "
$ws.Range("D48").Value = "'76.21"
$ws.Range("E48").Value = "  +9.60%  "
$ws.Range("D49").Value = "'106.23"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").Value = "2.822.49"
$ws.Range("E50").Value = "  -0.50%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'1.67"
$ws.Range("E51").Value = "  +0.69%  "
